$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.946.41"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.517.22"
$ws.Range("E3").Value = "  +3.41%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "534.18"
$ws.Range("E5").Value = "  +5.61%  "
$ws.Range("D6").Value = "133.98"
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D9").Value = "2.515.75"
$ws.Range("E9").Value = "  +2.75%  "
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "2.938.24"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "58.831.80"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "2.515.57"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").Value = "320.69"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +9.32%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "65.58"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").Value = "0.410"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "7.48"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("E29").Value = "  +5.25%  "
$ws.Range("D30").Value = "172.13"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +5.39%  "
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "18.13"
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "3.94"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +4.27%  "
$ws.Range("D40").Value = "36.71"
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +6.42%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.48"
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "276.77"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "5.10"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("D45").Value = "131.58"
$ws.Range("E45").Value = "  +10.06%  "
$ws.Range("D46").Value = "0.591"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "0.0934"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("E49").Value = "  +4.70%  "
$ws.Range("D50").Value = "17.00"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "1.753.62"
$ws.Range("E51").Value = "  +3.07%  "
